# Auto-generated Excel COM-interop script
# Refreshes the scraped cryptocurrency price (column D) and 1h volume-change
# percentage (column E) cells to the latest values from the data source,
# matching the GitHub Actions scheduled data-refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as plain text (the source feed already renders them
# with "."-grouped thousands, e.g. "26.911.58"); some of the new prices parse
# as valid numbers (e.g. "206.29"), so pin those specific cells to Text format
# first to keep them stored as strings, matching the rest of the column.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Price (column D) and Volume(1h) (column E) updates
$ws.Range("D2").Value = "26.911.58"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.549.45"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D5").Value = "206.29"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "0.487"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "22.10"
$ws.Range("E8").Value = "  +2.78%  "
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "1.770.48"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "1.549.61"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "26.902.59"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "61.62"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "217.17"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "9.23"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").Value = "154.26"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "6.63"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "14.94"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "3.21"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").Value = "1.420.41"
$ws.Range("E33").Value = "  +3.79%  "
$ws.Range("E34").Value = "  +4.73%  "
$ws.Range("E35").Value = "  +2.54%  "
$ws.Range("D36").Value = "0.967"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D39").Value = "0.525"
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").Value = "64.52"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").Value = "1.684.11"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").Value = "87.40"
$ws.Range("E48").Value = "  +1.49%  "
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("E50").Value = "  +3.31%  "
$ws.Range("E51").Value = "  +0.53%  "
